# The "Warmup Plan" sheet originally had two extra rows at the very top
# (a leftover "Properties"/"Value" pair) above the real "Phase" header row.
# This edit removes those two stray rows, which shifts the whole table
# (headers, data, merged cells, trailing spacer rows) up by two rows and
# drops the now-unused "Properties"/"Value" shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

# Delete rows 1 and 2 (the "Properties"/"Value" rows) - everything below
# shifts up by two rows automatically (headers, data, merged cells, and the
# trailing blank rows all renumber, so the sheet ends up one row 1..418
# instead of 1..420).
$ws.Rows("1:2").Delete()

# Leave the selection on the new header row, matching the state Excel is
# left in right after such a row deletion.
[void]$ws.Rows("1:1").Select()
